$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Rep ID" header back to "Address Name"
$ws.Range("A1").Value = "Address Name"
